$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.012.50"
$ws.Range("E2").Value = "  +1.93%  "
$ws.Range("D3").Value = "1.907.98"
$ws.Range("E3").Value = "  +2.27%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  -0.80%  "
$ws.Range("D5").Value = "'315.55"
$ws.Range("E5").Value = "  +1.22%  "
$ws.Range("E6").Value = "  -0.79%  "
$ws.Range("D7").Value = "'0.4810"
$ws.Range("E7").Value = "  +0.70%  "
$ws.Range("E8").Value = "  +1.34%  "
$ws.Range("D9").Value = "'0.07364"
$ws.Range("E9").Value = "  +0.53%  "
$ws.Range("D10").Value = "'0.9326"
$ws.Range("E10").Value = "  -0.26%  "
$ws.Range("E11").Value = "  +0.77%  "
$ws.Range("D12").Value = "'0.07752"
$ws.Range("E12").Value = "  -0.94%  "
$ws.Range("D13").Value = "1.906.83"
$ws.Range("E13").Value = "  +1.62%  "
$ws.Range("E14").Value = "  +1.15%  "
$ws.Range("D15").Value = "'6.636"
$ws.Range("E15").Value = "  +1.31%  "
$ws.Range("D16").Value = "'91.61"
$ws.Range("E16").Value = "  +1.47%  "
$ws.Range("D17").Value = "'1.006"
$ws.Range("E17").Value = "  -0.69%  "
$ws.Range("E18").Value = "  -0.58%  "
$ws.Range("D20").Value = "28.050.88"
$ws.Range("E20").Value = "  +1.75%  "
$ws.Range("E21").Value = "  +0.98%  "
$ws.Range("E22").Value = "  +0.96%  "
$ws.Range("D23").Value = "2.164.70"
$ws.Range("E23").Value = "  +2.95%  "
$ws.Range("D24").Value = "'10.89"
$ws.Range("E24").Value = "  +1.84%  "
$ws.Range("D25").Value = "'155.84"
$ws.Range("E25").Value = "  +0.87%  "
$ws.Range("D26").Value = "'1.921"
$ws.Range("E26").Value = "  -1.14%  "
$ws.Range("D27").Value = "'18.48"
$ws.Range("D28").Value = "'2.128"
$ws.Range("E28").Value = "  +5.25%  "
$ws.Range("D29").Value = "'116.90"
$ws.Range("E29").Value = "  +1.23%  "
$ws.Range("D30").Value = "'4.959"
$ws.Range("E30").Value = "  -0.38%  "
$ws.Range("D31").Value = "'0.08935"
$ws.Range("E31").Value = "  +0.44%  "
$ws.Range("D32").Value = "'3.296"
$ws.Range("E32").Value = "  -0.96%  "
$ws.Range("D33").Value = "'1.258"
$ws.Range("E33").Value = "  +3.32%  "
$ws.Range("D34").Value = "'0.7734"
$ws.Range("E34").Value = "  +2.06%  "
$ws.Range("D35").Value = "'4.679"
$ws.Range("E35").Value = "  +1.46%  "
$ws.Range("D36").Value = "'2.641"
$ws.Range("E36").Value = "  -3.98%  "
$ws.Range("D37").Value = "'0.02062"
$ws.Range("E37").Value = "  +1.48%  "
$ws.Range("D38").Value = "'1.111"
$ws.Range("E38").Value = "  -0.71%  "
$ws.Range("E39").Value = "  +0.90%  "
$ws.Range("D40").Value = "'0.5491"
$ws.Range("E40").Value = "  +3.43%  "
$ws.Range("D41").Value = "'2.990"
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").Value = "'7.031"
$ws.Range("E42").Value = "  -0.62%  "
$ws.Range("E43").Value = "  +0.33%  "
$ws.Range("D44").Value = "'8.491"
$ws.Range("E44").Value = "  +0.23%  "
$ws.Range("D45").Value = "'10.70"
$ws.Range("E45").Value = "  +1.08%  "
$ws.Range("D46").Value = "'0.4825"
$ws.Range("E46").Value = "  +0.52%  "
$ws.Range("D47").Value = "'107.96"
$ws.Range("E47").Value = "  +4.92%  "
$ws.Range("E48").Value = "  -0.73%  "
$ws.Range("D49").Value = "'1.648"
$ws.Range("E49").Value = "  -0.21%  "
$ws.Range("D50").Value = "'67.89"
$ws.Range("E50").Value = "  +0.68%  "
$ws.Range("D51").Value = "'0.06073"
$ws.Range("E51").Value = "  -0.06%  "
